# Generate Report for Handoff
#
# The "c6bb5ffc-54de-40a4-8f3d-268c64ed5a66.md" file (row 3 in every sheet)
# moves from "In Translation" to "Ready for handoff", and the per-language
# sheets record a fresh "Latest Handoff Datetime" for that row.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: update the zh-cn / de-de status columns for row 3 ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet: status + new handoff datetime for row 3 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-09 07:37:05"

# --- de-de sheet: status + new handoff datetime for row 3 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-09 07:37:09"
